# Apply the 2021-03-17 -> 2021-03-18 "as of" update to the disclaimer text,
# and refresh the Weight (D) / Percent Change (E) columns of the holdings
# table on Sheet1 with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (password hash "D382"); unprotect so the cells
# can be edited, then restore protection afterwards.
$ws.Unprotect("D382")

# --- Disclaimer text (A16): bump the "as of" date ---------------------------
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-18 for illustrative purposes only and are subject to change."

# --- Weight (D) / Percent Change (E) updates --------------------------------
$ws.Range("D2").Value = 0.03033676084725239
$ws.Range("E2").Value = -0.010435408420295

$ws.Range("D3").Value = 0.02335343428455921
$ws.Range("E3").Value = -0.005183905208590467

$ws.Range("D4").Value = 0.05094465148865769
$ws.Range("E4").Value = -0.006385068762278845

$ws.Range("D5").Value = 0.1362642845788367
$ws.Range("E5").Value = -0.009186798230690685

$ws.Range("D6").Value = 0.03227203548992941
$ws.Range("E6").Value = -0.04975124378109452

$ws.Range("D7").Value = 0.1208025710770028
$ws.Range("E7").Value = -0.02450032237266275

$ws.Range("D8").Value = 0.1004063855967878
$ws.Range("E8").Value = -0.004632310364794501

$ws.Range("D9").Value = 0.02769927760213363
$ws.Range("E9").Value = -0.006323396567299078

$ws.Range("D10").Value = 0.1222149609559751
$ws.Range("E10").Value = 0.002207948615013944

$ws.Range("D11").Value = 0.2500365516442843
$ws.Range("E11").Value = -0.03141019666324929

$ws.Range("D12").Value = 0.1056690864345809
$ws.Range("E12").Value = -0.02125585023400944

$ws.Range("E13").Value = -0.01705023814990791

# Restore sheet protection.
$ws.Protect("D382")
